$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2: Investigate web technologies
$t.Cell(2, 1).Range.Text = "Investigate web technologies"
$t.Cell(2, 2).Range.Text = "07/05/2025"
$t.Cell(2, 3).Range.Text = "14/05/2025"
$t.Cell(2, 4).Range.Text = "No"

# Row 3: Web interface design
$t.Cell(3, 1).Range.Text = "Web interface design"
$t.Cell(3, 2).Range.Text = "07/05/2025"
$t.Cell(3, 3).Range.Text = "21/05/2025"
$t.Cell(3, 4).Range.Text = "No"

# Row 4: Review web design
$t.Cell(4, 1).Range.Text = "Review web design"
$t.Cell(4, 2).Range.Text = "14/05/2025"
$t.Cell(4, 3).Range.Text = "21/05/2025"
$t.Cell(4, 4).Range.Text = "Yes"

# Row 5: Create basic web pages
$t.Cell(5, 1).Range.Text = "Create basic web pages"
$t.Cell(5, 2).Range.Text = "21/05/2025"
$t.Cell(5, 3).Range.Text = "28/05/2025"
$t.Cell(5, 4).Range.Text = "No"

# Row 6: Add styles formatting
$t.Cell(6, 1).Range.Text = "Add styles formatting"
$t.Cell(6, 2).Range.Text = "21/05/2025"
$t.Cell(6, 3).Range.Text = "04/06/2025"
$t.Cell(6, 4).Range.Text = "No"

# Row 7: Present draft to client for feedback
$t.Cell(7, 1).Range.Text = "Present draft to client for feedback"
$t.Cell(7, 2).Range.Text = "21/05/2025"
$t.Cell(7, 3).Range.Text = "04/06/2025"
$t.Cell(7, 4).Range.Text = "Yes"

# Row 8: Create & update BlueSpringsHotel_ProjectPlanV2.docx (mixed formatting)
$t.Cell(8, 1).Range.Text = "Create & update "
$t.Cell(8, 1).Range.InsertAfter("BlueSpringsHotel_ProjectPlan")
$t.Cell(8, 1).Range.InsertAfter("V")
$t.Cell(8, 1).Range.InsertAfter("2.docx")
$t.Cell(8, 2).Range.Text = "04/06/2025"
$t.Cell(8, 3).Range.Text = "04/06/2025"
$t.Cell(8, 4).Range.Text = "No"

# Apply bold to "BlueSpringsHotel_ProjectPlan", "V" and "2.docx" as separate runs
$cellRange = $t.Cell(8, 1).Range
$cellStart = $cellRange.Start
$prefixLen = "Create & update ".Length
$nameLen = "BlueSpringsHotel_ProjectPlan".Length
$vLen = "V".Length
$suffixLen = "2.docx".Length

$r1 = $d.Range($cellStart + $prefixLen, $cellStart + $prefixLen + $nameLen)
$r1.Font.Bold = 1

$r2 = $d.Range($cellStart + $prefixLen + $nameLen, $cellStart + $prefixLen + $nameLen + $vLen)
$r2.Font.Bold = 1

$r3 = $d.Range($cellStart + $prefixLen + $nameLen + $vLen, $cellStart + $prefixLen + $nameLen + $vLen + $suffixLen)
$r3.Font.Bold = 1
